# Apply the edits described by the diff:
#  - MultC_2 (sheet2): clear C6 (was "Wow, that's way off.")
#  - Matching (sheet5): A1 text "Key terms" -> "Terms"; selection moves to B8
#  - MultC_2 becomes the active/selected sheet (was Matching)

$wb = $excel.ActiveWorkbook

$sheetMultC2 = $wb.Worksheets.Item("MultC_2")
$sheetMatching = $wb.Worksheets.Item("Matching")

# Clear the stray comment cell on MultC_2
$sheetMultC2.Range("C6").ClearContents()

# Rename "Key terms" header to "Terms" on Matching
$sheetMatching.Range("A1").Value = "Terms"

# Update the lingering selection on Matching to B8
$sheetMatching.Range("B8").Select()

# Make MultC_2 the active sheet/tab and put the selection on its last used cell
$sheetMultC2.Activate()
$sheetMultC2.Range("C6").Select()
